$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'26.388.44"
$ws.Range("D2").Style = $style
$style = $ws.Range("E2").Style
$ws.Range("E2").Value = "'  -0.44%  "
$ws.Range("E2").Style = $style
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.723.01"
$ws.Range("D3").Style = $style
$style = $ws.Range("E3").Style
$ws.Range("E3").Value = "'  -0.45%  "
$ws.Range("E3").Style = $style
$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = $style
$style = $ws.Range("E4").Style
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = $style
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'242.21"
$ws.Range("D5").Style = $style
$style = $ws.Range("E5").Style
$ws.Range("E5").Value = "'  -2.02%  "
$ws.Range("E5").Style = $style
$style = $ws.Range("E6").Style
$ws.Range("E6").Value = "'  +0.03%  "
$ws.Range("E6").Style = $style
$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.4857"
$ws.Range("D7").Style = $style
$style = $ws.Range("E7").Style
$ws.Range("E7").Value = "'  +0.24%  "
$ws.Range("E7").Style = $style
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.2584"
$ws.Range("D8").Style = $style
$style = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  -3.20%  "
$ws.Range("E8").Style = $style
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.06190"
$ws.Range("D9").Style = $style
$style = $ws.Range("E9").Style
$ws.Range("E9").Value = "'  -0.50%  "
$ws.Range("E9").Style = $style
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'1.730.17"
$ws.Range("D10").Style = $style
$style = $ws.Range("E10").Style
$ws.Range("E10").Value = "'  +0.00%  "
$ws.Range("E10").Style = $style
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.06969"
$ws.Range("D11").Style = $style
$style = $ws.Range("E11").Style
$ws.Range("E11").Value = "'  -1.33%  "
$ws.Range("E11").Style = $style
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'15.46"
$ws.Range("D12").Style = $style
$style = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  -1.28%  "
$ws.Range("E12").Style = $style
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'4.529"
$ws.Range("D13").Style = $style
$style = $ws.Range("E13").Style
$ws.Range("E13").Value = "'  -1.94%  "
$ws.Range("E13").Style = $style
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.5964"
$ws.Range("D14").Style = $style
$style = $ws.Range("E14").Style
$ws.Range("E14").Value = "'  -2.53%  "
$ws.Range("E14").Style = $style
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'77.16"
$ws.Range("D15").Style = $style
$style = $ws.Range("E15").Style
$ws.Range("E15").Value = "'  -0.28%  "
$ws.Range("E15").Style = $style
$style = $ws.Range("E16").Style
$ws.Range("E16").Value = "'  -0.01%  "
$ws.Range("E16").Style = $style
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'26.387.86"
$ws.Range("D17").Style = $style
$style = $ws.Range("E17").Style
$ws.Range("E17").Value = "'  -0.43%  "
$ws.Range("E17").Style = $style
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = $style
$style = $ws.Range("E18").Style
$ws.Range("E18").Value = "'  +0.06%  "
$ws.Range("E18").Style = $style
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.000007206"
$ws.Range("D19").Style = $style
$style = $ws.Range("E19").Style
$ws.Range("E19").Value = "'  -0.31%  "
$ws.Range("E19").Style = $style
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'11.31"
$ws.Range("D20").Style = $style
$style = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  -2.15%  "
$ws.Range("E20").Style = $style
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'1.938.58"
$ws.Range("D21").Style = $style
$style = $ws.Range("E21").Style
$ws.Range("E21").Value = "'  -0.72%  "
$ws.Range("E21").Style = $style
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.439"
$ws.Range("D22").Style = $style
$style = $ws.Range("E22").Style
$ws.Range("E22").Value = "'  -1.63%  "
$ws.Range("E22").Style = $style
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'8.468"
$ws.Range("D23").Style = $style
$style = $ws.Range("E23").Style
$ws.Range("E23").Value = "'  -3.68%  "
$ws.Range("E23").Style = $style
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'5.097"
$ws.Range("D24").Style = $style
$style = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  -3.06%  "
$ws.Range("E24").Style = $style
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'137.69"
$ws.Range("D25").Style = $style
$style = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  +0.17%  "
$ws.Range("E25").Style = $style
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'15.23"
$ws.Range("D26").Style = $style
$style = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  -1.32%  "
$ws.Range("E26").Style = $style
$style = $ws.Range("E27").Style
$ws.Range("E27").Value = "'  +0.24%  "
$ws.Range("E27").Style = $style
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'106.55"
$ws.Range("D28").Style = $style
$style = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  -1.59%  "
$ws.Range("E28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'1.721"
$ws.Range("D29").Style = $style
$style = $ws.Range("E29").Style
$ws.Range("E29").Value = "'  -3.26%  "
$ws.Range("E29").Style = $style
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'3.919"
$ws.Range("D30").Style = $style
$style = $ws.Range("E30").Style
$ws.Range("E30").Value = "'  -1.64%  "
$ws.Range("E30").Style = $style
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'0.07996"
$ws.Range("D31").Style = $style
$style = $ws.Range("E31").Style
$ws.Range("E31").Value = "'  +0.07%  "
$ws.Range("E31").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.659"
$ws.Range("D32").Style = $style
$style = $ws.Range("E32").Style
$ws.Range("E32").Value = "'  -0.92%  "
$ws.Range("E32").Style = $style
$style = $ws.Range("E33").Style
$ws.Range("E33").Value = "'  -1.48%  "
$ws.Range("E33").Style = $style
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'2.606"
$ws.Range("D34").Style = $style
$style = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  -0.23%  "
$ws.Range("E34").Style = $style
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.9966"
$ws.Range("D35").Style = $style
$style = $ws.Range("E35").Style
$ws.Range("E35").Value = "'  -0.89%  "
$ws.Range("E35").Style = $style
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.6230"
$ws.Range("D36").Style = $style
$style = $ws.Range("E36").Style
$ws.Range("E36").Value = "'  -1.87%  "
$ws.Range("E36").Style = $style
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.9353"
$ws.Range("D37").Style = $style
$style = $ws.Range("E37").Style
$ws.Range("E37").Value = "'  +4.29%  "
$ws.Range("E37").Style = $style
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'1.951"
$ws.Range("D38").Style = $style
$style = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  -3.22%  "
$ws.Range("E38").Style = $style
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'2.387"
$ws.Range("D39").Style = $style
$style = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  -0.04%  "
$ws.Range("E39").Style = $style
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.9995"
$ws.Range("D40").Style = $style
$style = $ws.Range("E40").Style
$ws.Range("E40").Value = "'  -0.40%  "
$ws.Range("E40").Style = $style
$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.01470"
$ws.Range("D41").Style = $style
$style = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  -2.25%  "
$ws.Range("E41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'99.89"
$ws.Range("D42").Style = $style
$style = $ws.Range("E42").Style
$ws.Range("E42").Value = "'  -1.56%  "
$ws.Range("E42").Style = $style
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'5.331"
$ws.Range("D43").Style = $style
$style = $ws.Range("E43").Style
$ws.Range("E43").Value = "'  -2.64%  "
$ws.Range("E43").Style = $style
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.3829"
$ws.Range("D44").Style = $style
$style = $ws.Range("E44").Style
$ws.Range("E44").Value = "'  -1.75%  "
$ws.Range("E44").Style = $style
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'6.839"
$ws.Range("D45").Style = $style
$style = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  -2.64%  "
$ws.Range("E45").Style = $style
$style = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  -1.73%  "
$ws.Range("E46").Style = $style
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.05364"
$ws.Range("D47").Style = $style
$style = $ws.Range("E47").Style
$ws.Range("E47").Value = "'  -0.34%  "
$ws.Range("E47").Style = $style
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'7.710"
$ws.Range("D48").Style = $style
$style = $ws.Range("E48").Style
$ws.Range("E48").Value = "'  -2.81%  "
$ws.Range("E48").Style = $style
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'30.06"
$ws.Range("D49").Style = $style
$style = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  -1.69%  "
$ws.Range("E49").Style = $style
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'1.225"
$ws.Range("D50").Style = $style
$style = $ws.Range("E50").Style
$ws.Range("E50").Value = "'  -2.07%  "
$ws.Range("E50").Style = $style
$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'50.79"
$ws.Range("D51").Style = $style
$style = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  -1.48%  "
$ws.Range("E51").Style = $style
